# Update "Lương" worksheet: insert "Ứng lương" rows after each location's
# "Công phụ phẫu 2" row, and append "Tổng lương" summary rows at the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

# Insert a new row before row 11 ("Lương cơ bản tại LONG XUYÊN" and below
# shift down by one) and fill it with the new "Ứng lương tại CẦN THƠ" entry.
$ws.Rows.Item(11).Insert()
$ws.Cells.Item(11, 1).Value = "Ứng lương tại CẦN THƠ"

# Insert a new row before the ("Lương cơ bản tại SÓC TRĂNG") row, which is
# now row 19 after the previous insertion, for "Ứng lương tại LONG XUYÊN".
$ws.Rows.Item(19).Insert()
$ws.Cells.Item(19, 1).Value = "Ứng lương tại LONG XUYÊN"

# Append "Ứng lương tại SÓC TRĂNG" right after the last existing row, which
# is now row 26 ("Công phụ phẫu 2 tại SÓC TRĂNG") after the two insertions.
$ws.Cells.Item(27, 1).Value = "Ứng lương tại SÓC TRĂNG"

# Append the four new "Tổng lương" summary rows at the bottom.
$ws.Cells.Item(28, 1).Value = "Tổng lương tại CẦN THƠ"
$ws.Cells.Item(29, 1).Value = "Tổng lương tại LONG XUYÊN"
$ws.Cells.Item(30, 1).Value = "Tổng lương tại SÓC TRĂNG"
$ws.Cells.Item(31, 1).Value = "Tổng lương"
